# Tweak the Snowflake + SageMaker Autopilot serverless architecture diagram:
# use the short names for SageMaker and API Gateway so they match the
# other labels on the diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    if (-not $sh.TextFrame.HasText) { continue }

    $tr = $sh.TextFrame.TextRange
    $text = $tr.Text

    if ($sh.Name -eq "TextBox 15" -and $text -eq "Amazon SageMaker") {
        # Whole textbox is just this one run -- replace it outright.
        $tr.Text = "SageMaker"
    }
    elseif ($sh.Name -eq "TextBox 16" -and $text.StartsWith("Amazon API Gateway")) {
        # Textbox also contains a line break + "endpoint" after the title;
        # only touch the "Amazon API Gateway" prefix so the rest is untouched.
        $prefixLen = "Amazon API Gateway".Length
        $sub = $tr.Characters(1, $prefixLen)
        $sub.Text = "API Gateway"
    }
}
